# Cronograma de actividades - actualizacion de cronograma de desarrollo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths / visibility -------------------------------------------------
# (ColumnWidth is quantized to whole pixels by this host, so we feed it the
#  "character width minus the 5-pixel padding" value that lands closest to the
#  canonical OOXML width stored in the target file.)
$ws.Columns("A").ColumnWidth = 4.451822916666667
$ws.Columns("B").Hidden = $true
$ws.Columns("C").ColumnWidth = 33.736979166666664
$ws.Columns("D").ColumnWidth = 66.02213541666667

# --- Row 17 -----------------------------------------------------------------
$c17 = $ws.Range("C17")
$c17.Value = "Actualización de plantilla administrativa y sincronizacion con los controladores de enrutamiento para el modulo multimedia"
$c17.Characters(113, 10).Font.Bold = $true
$c17.Borders.Item(7).LineStyle = 1
$c17.Borders.Item(10).LineStyle = 1
$c17.HorizontalAlignment = -4108
$c17.VerticalAlignment = -4108
$c17.WrapText = $true

$ws.Range("G16").Copy() | Out-Null
$ws.Range("G17").PasteSpecial(-4122) | Out-Null
$ws.Range("G17").Value = "Terminado"

$ws.Range("H17").Value = "Pendiente para cambiar colores"

$ws.Rows(17).RowHeight = 60

# --- Row 18 -----------------------------------------------------------------
$c18 = $ws.Range("C18")
$c18.Value = "Creacion de controladores del CRUD (CREAR CONSULTAR ELIMINAR) para el modulo multimedia"
$c18.Characters(31, 32).Font.Bold = $true
$c18.Characters(78, 10).Font.Bold = $true
$c18.Borders.Item(7).LineStyle = 1
$c18.Borders.Item(10).LineStyle = 1
$c18.HorizontalAlignment = -4108
$c18.VerticalAlignment = -4108
$c18.WrapText = $true

$d18 = $ws.Range("D18")
$d18.Value = "Dento de esta creacion de controladores se crearon sus respectivos modelos y formularios, se configuraron los archivos js que hacen el respectivo redireccionamiento al controlador correspondiente. para la creacion de elementos multimedia, ya se creaN los elementos DESDE FORMULARIOS ESTABLECIDOS CON SUS RESPECTIVAS VALIDACIONES y son sincronizados con la vista que muestra los archivos dentro del panel de administración"
$d18.HorizontalAlignment = -4108
$d18.VerticalAlignment = -4108
$d18.WrapText = $true

$ws.Range("E11").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = 44102

$ws.Range("F11").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122) | Out-Null
$ws.Range("F18").Value = 44102

$ws.Range("G16").Copy() | Out-Null
$ws.Range("G18").PasteSpecial(-4122) | Out-Null
$ws.Range("G18").Value = "Terminado"

$ws.Rows(18).RowHeight = 114

# --- View state ---------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("H21").Select()
